# Update Excel file with latest predictions
# Rewrites the data rows of all five sheets (Home win, Draw, Btts, Over_Under,
# Away Win) to the latest betclever prediction feed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Home win"  (old last row 5 -> new last row 7: +2 rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")
$ws.Rows("2:3").Insert()

$data1 = @(
    @("15-01-2025 22:30", "BRAZIL", "CATARINENSE - 1", "Caravaggio - Chapecoense-sc", 70, 3.3),
    @("15-01-2025 14:00", "BRAZIL", "PARANAENSE - 1", "Andraus Brasil - Azuriz", 73.3, 3),
    @("15-01-2025 11:30", "IRAQ", "IRAQI LEAGUE", "Al Talaba - Erbil", 78.3, 2),
    @("15-01-2025 12:35", "SAUDI-ARABIA", "DIVISION 1", "Al Najma - Al Jandal", 86.7, 1.77),
    @("15-01-2025 18:30", "SPAIN", "COPA DEL REY", "Pontevedra - Getafe", 70, 3.5),
    @("15-01-2025 12:00", "THAILAND", "THAI LEAGUE 1", "Bangkok Glass - Ratchaburi", 73.3, 1.7)
)

$r = 2
foreach ($row in $data1) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 2: "Draw"  (old last row 6 -> new last row 4: -2 rows)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")
$ws.Rows("5:6").Delete()

$data2 = @(
    @("15-01-2025 16:00", "ALBANIA", "CUP", "Teuta Durrës - AF Elbasani", 80, 2.8),
    @("15-01-2025 23:30", "BRAZIL", "CATARINENSE - 1", "Avai - Santa Catarina", 70, 4.2),
    @("16-01-2025 19:30", "ENGLAND", "PREMIER LEAGUE", "Ipswich - Brighton", 60, 3.7)
)

$r = 2
foreach ($row in $data2) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 3: "Btts"  (old last row 10 -> new last row 9: -1 row)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")
$ws.Rows("10:10").Delete()

$data3 = @(
    @("15-01-2025 19:30", "GERMANY", "BUNDESLIGA", "Union Berlin - FC Augsburg", 80, 1.95),
    @("15-01-2025 21:30", "BRAZIL", "PAULISTA - A1", "Velo Clube - Noroeste", 81.7, 2.05),
    @("15-01-2025 23:00", "COSTA-RICA", "PRIMERA DIVISIÓN", "Santa Ana - Puntarenas FC", 83.3, 1.83),
    @("15-01-2025 17:30", "FRANCE", "COUPE DE FRANCE", "Thaon - Strasbourg", 84, 2.6),
    @("15-01-2025 20:15", "PORTUGAL", "TAÇA DE PORTUGAL", "SC Braga - Lusitano Évora 1911", 84, 2.4),
    @("15-01-2025 17:00", "SAUDI-ARABIA", "PRO LEAGUE", "Al Kholood - Al-Ahli Jeddah", 76.7, 1.85),
    @("15-01-2025 11:00", "THAILAND", "THAI LEAGUE 1", "Port FC - Khon Kaen United", 87.8, 1.75),
    @("16-01-2025 13:55", "SAUDI-ARABIA", "PRO LEAGUE", "Al Akhdoud - Al-Fayha", 76.7, 1.85)
)

$r = 2
foreach ($row in $data3) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 4: "Over_Under"  (old last row 7 -> new last row 6: -1 row)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")
$ws.Rows("7:7").Delete()

$data4 = @(
    @("15-01-2025 17:45", "NETHERLANDS", "KNVB BEKER", "GO Ahead Eagles - Twente", 80, 1.73, 50, 2.8),
    @("15-01-2025 17:30", "FRANCE", "COUPE DE FRANCE", "Cannes - Lorient", 80, 1.8, 60, 3),
    @("15-01-2025 15:00", "PORTUGAL", "LIGA 3", "SC Covilha - Lusitânia", 93.3, 1.95, 66.7, $null),
    @("15-01-2025 12:35", "SAUDI-ARABIA", "DIVISION 1", "Al-Hazm - Al Safa", 86.7, 1.8, 33.3, 3),
    @("15-01-2025 11:00", "WORLD", "FRIENDLIES CLUBS", "Wisla Plock - Warta Gorzów", 73.3, 1.6, 60, 2.6)
)

$r = 2
foreach ($row in $data4) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    if ($row[7] -eq $null) {
        $ws.Range("H$r").Value = ""
    } else {
        $ws.Range("H$r").Value = $row[7]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 5: "Away Win"  (row count unchanged: old/new last row 2)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Away Win")

$ws.Range("A2").Value = "16-01-2025 14:00"
$ws.Range("B2").Value = "INDIA"
$ws.Range("C2").Value = "INDIAN SUPER LEAGUE"
$ws.Range("D2").Value = "Minerva Punjab - Mumbai City"
$ws.Range("E2").Value = 73.3
$ws.Range("F2").Value = 2.05
